$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text updates (reordering of countries / timestamp refresh) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 20:05"

# Swap Portugal <-> Emiratos Arabes Unidos (EAU overtakes Portugal in ranking)
$ws.Range("A30").Value = "Emiratos Arabes Unidos"
$ws.Range("A31").Value = "Portugal"

# Sudan del Sur moves up right after Guinea Ecuatorial; rows below shift down by one
$ws.Range("A115").Value = "Sudan del Sur"
$ws.Range("A116").Value = "Costa Rica"
$ws.Range("A117").Value = "Niger"
$ws.Range("A118").Value = "Republica de Chipre"
$ws.Range("A119").Value = "Nepal"
$ws.Range("A120").Value = "Paraguay"
$ws.Range("A121").Value = "Burkina Faso"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 1735764 ; $ws.Range("C4").Value = 10489 ; $ws.Range("D4").Value = 482901 ; $ws.Range("E4").Value = 1151471 ; $ws.Range("G4").Value = 820 ; $ws.Range("H4").Value = 101392
$ws.Range("B11").Value = 181719 ; $ws.Range("C11").Value = 431 ; $ws.Range("E11").Value = 10399 ; $ws.Range("G11").Value = 22 ; $ws.Range("H11").Value = 8520
$ws.Range("B13").Value = 158042 ; $ws.Range("C13").Value = 7249 ; $ws.Range("D13").Value = 67711 ; $ws.Range("E13").Value = 85797 ; $ws.Range("G13").Value = 190 ; $ws.Range("H13").Value = 4534
$ws.Range("B30").Value = 31969 ; $ws.Range("C30").Value = 883 ; $ws.Range("D30").Value = 16371 ; $ws.Range("E30").Value = 15343 ; $ws.Range("G30").Value = 2 ; $ws.Range("H30").Value = 255
$ws.Range("B31").Value = 31292 ; $ws.Range("C31").Value = 285 ; $ws.Range("D31").Value = 18349 ; $ws.Range("E31").Value = 11587 ; $ws.Range("G31").Value = 14 ; $ws.Range("H31").Value = 1356
$ws.Range("B33").Value = 24803 ; $ws.Range("C33").Value = 68 ; $ws.Range("E33").Value = 2112 ; $ws.Range("G33").Value = 16 ; $ws.Range("H33").Value = 1631
$ws.Range("E41").Value = 5205 ; $ws.Range("G41").Value = 11 ; $ws.Range("H41").Value = 1227
$ws.Range("B115").Value = 994 ; $ws.Range("C115").Value = 188 ; $ws.Range("D115").Value = 6 ; $ws.Range("E115").Value = 978 ; $ws.Range("G115").Value = 2
$ws.Range("B116").Value = 956 ; $ws.Range("D116").Value = 634 ; $ws.Range("E116").Value = 312 ; $ws.Range("H116").Value = 10
$ws.Range("B117").Value = 952 ; $ws.Range("D117").Value = 796 ; $ws.Range("E117").Value = 93 ; $ws.Range("H117").Value = 63
$ws.Range("B118").Value = 939 ; $ws.Range("C118").Value = 0 ; $ws.Range("D118").Value = 594 ; $ws.Range("E118").Value = 328 ; $ws.Range("H118").Value = 17
$ws.Range("B119").Value = 886 ; $ws.Range("C119").Value = 114 ; $ws.Range("D119").Value = 183 ; $ws.Range("E119").Value = 699 ; $ws.Range("H119").Value = 4
$ws.Range("B120").Value = 884 ; $ws.Range("C120").Value = 7 ; $ws.Range("D120").Value = 392 ; $ws.Range("E120").Value = 481 ; $ws.Range("G120").Value = 0 ; $ws.Range("H120").Value = 11
$ws.Range("B121").Value = 845 ; $ws.Range("C121").Value = 13 ; $ws.Range("D121").Value = 672 ; $ws.Range("E121").Value = 120 ; $ws.Range("G121").Value = 1 ; $ws.Range("H121").Value = 53
$ws.Range("B145").Value = 346 ; $ws.Range("C145").Value = 7 ; $ws.Range("D145").Value = 245 ; $ws.Range("E145").Value = 101
$ws.Range("B177").Value = 79 ; $ws.Range("C177").Value = 2 ; $ws.Range("E177").Value = 36
